$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.209.05"
$ws.Range("E2").Value = "  -0.01%  "
$ws.Range("D3").Value = "3.030.12"
$ws.Range("E3").Value = "  +0.94%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'576.67"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.89%  "
$ws.Range("D6").Value = "'168.18"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +3.59%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").Value = "3.027.11"
$ws.Range("E8").Value = "  +0.90%  "
$ws.Range("D9").Value = "'0.519"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.17%  "
$ws.Range("D10").Value = "'6.68"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.01%  "
$ws.Range("E11").Value = "  -1.34%  "
$ws.Range("E12").Value = "  +5.58%  "
$ws.Range("E13").Value = "  -2.08%  "
$ws.Range("D14").Value = "'36.28"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +5.06%  "
$ws.Range("E15").Value = "  -0.45%  "
$ws.Range("D16").Value = "66.169.50"
$ws.Range("E16").Value = "  +0.00%  "
$ws.Range("D17").Value = "3.530.14"
$ws.Range("E17").Value = "  +0.83%  "
$ws.Range("D18").Value = "'7.26"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +5.13%  "
$ws.Range("D19").Value = "'16.51"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +19.44%  "
$ws.Range("D20").Value = "3.029.23"
$ws.Range("E20").Value = "  +0.84%  "
$ws.Range("D21").Value = "'474.04"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +4.33%  "
$ws.Range("D22").Value = "'0.707"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +3.00%  "
$ws.Range("E23").Value = "  +2.02%  "
$ws.Range("D24").Value = "'83.22"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.27%  "
$ws.Range("D25").Value = "'12.81"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +4.73%  "
$ws.Range("E26").Value = "  -0.52%  "
$ws.Range("D27").Value = "'10.06"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -3.66%  "
$ws.Range("E28").Value = "  +0.09%  "
$ws.Range("D29").Value = "'8.20"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.45%  "
$ws.Range("E30").Value = "  +2.00%  "
$ws.Range("E31").Value = "  +0.99%  "
$ws.Range("D32").Value = "'0.118"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +7.03%  "
$ws.Range("D33").Value = "0.0₃0997"
$ws.Range("E33").Value = "  -5.87%  "
$ws.Range("D34").Value = "'28.01"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +3.15%  "
$ws.Range("D35").Value = "'0.999"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("D36").Value = "'0.993"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.04%  "
$ws.Range("D37").Value = "'5.86"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +1.13%  "
$ws.Range("D38").Value = "'48.10"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +9.46%  "
$ws.Range("D39").Value = "'2.05"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -4.69%  "
$ws.Range("E40").Value = "  -0.48%  "
$ws.Range("D41").Value = "'0.310"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.17%  "
$ws.Range("E42").Value = "  -0.69%  "
$ws.Range("E43").Value = "  -5.03%  "
$ws.Range("D44").Value = "'8.61"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +2.44%  "
$ws.Range("E45").Value = "  +0.12%  "
$ws.Range("D46").Value = "'382.38"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -3.58%  "
$ws.Range("D47").Value = "2.720.36"
$ws.Range("E47").Value = "  -2.54%  "
$ws.Range("D48").Value = "'134.39"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.27%  "
$ws.Range("E49").Value = "  +0.01%  "
$ws.Range("D50").Value = "'24.51"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +2.90%  "
$ws.Range("D51").Value = "'2.23"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +4.17%  "
